$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.4
$ws.Range("W2").Value = 7.5
$ws.Range("Y2").Value = 9.5
$ws.Range("AJ2").Value = 12
$ws.Range("AK2").Value = 34
$ws.Range("AW2").Value = 5
